$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 103
    $ws.Cells.Item(103, 2).Value = 7630582
    $ws.Cells.Item(103, 6).Value = "Athlone Town"
    $ws.Cells.Item(103, 7).Value = "Wexford FC"
    $ws.Cells.Item(103, 8).Value = 3
    $ws.Cells.Item(103, 9).Value = 0
    $ws.Cells.Item(103, 10).Value = "H"
    $ws.Cells.Item(103, 11).Value = 1.95
    $ws.Cells.Item(103, 12).Value = 3.5
    $ws.Cells.Item(103, 13).Value = 3.4
    $ws.Cells.Item(103, 14).Value = 3
    $ws.Cells.Item(103, 15).Value = 3.3
    $ws.Cells.Item(103, 16).Value = 2.15
    $ws.Cells.Item(103, 17).Value = 0.25
    $ws.Cells.Item(103, 18).Value = 1.925
    $ws.Cells.Item(103, 19).Value = 1.925
    $ws.Cells.Item(103, 20).Value = 2.75
    $ws.Cells.Item(103, 21).Value = 2
    $ws.Cells.Item(103, 22).Value = 1.85
    $ws.Cells.Item(103, 23).Value = 2
    $ws.Cells.Item(103, 24).Value = -1
    $ws.Cells.Item(103, 25).Value = -1
    $ws.Cells.Item(103, 26).Value = 0.925
    $ws.Cells.Item(103, 27).Value = -1
    $ws.Cells.Item(103, 28).Value = 0.5
    $ws.Cells.Item(103, 29).Value = -0.5
    # Row 104
    $ws.Cells.Item(104, 2).Value = 7630583
    $ws.Cells.Item(104, 6).Value = "Cork City"
    $ws.Cells.Item(104, 7).Value = "Kerry FC"
    $ws.Cells.Item(104, 8).Value = 2
    $ws.Cells.Item(104, 9).Value = 0
    $ws.Cells.Item(104, 10).Value = "H"
    $ws.Cells.Item(104, 11).Value = 1.125
    $ws.Cells.Item(104, 12).Value = 8
    $ws.Cells.Item(104, 13).Value = 15
    $ws.Cells.Item(104, 14).Value = 1.285
    $ws.Cells.Item(104, 15).Value = 5.5
    $ws.Cells.Item(104, 16).Value = 9
    $ws.Cells.Item(104, 17).Value = -1.5
    $ws.Cells.Item(104, 18).Value = 1.8
    $ws.Cells.Item(104, 19).Value = 2
    $ws.Cells.Item(104, 20).Value = 2.75
    $ws.Cells.Item(104, 21).Value = 1.775
    $ws.Cells.Item(104, 22).Value = 2.025
    $ws.Cells.Item(104, 23).Value = 0.2849999999999999
    $ws.Cells.Item(104, 24).Value = -1
    $ws.Cells.Item(104, 25).Value = -1
    $ws.Cells.Item(104, 26).Value = 0.8
    $ws.Cells.Item(104, 27).Value = -1
    $ws.Cells.Item(104, 28).Value = -1
    $ws.Cells.Item(104, 29).Value = 1.025
    # Row 105
    $ws.Cells.Item(105, 2).Value = 7630584
    $ws.Cells.Item(105, 6).Value = "Treaty United FC"
    $ws.Cells.Item(105, 7).Value = "Cobh Ramblers"
    $ws.Cells.Item(105, 8).Value = 3
    $ws.Cells.Item(105, 9).Value = 1
    $ws.Cells.Item(105, 10).Value = "H"
    $ws.Cells.Item(105, 11).Value = 2.5
    $ws.Cells.Item(105, 12).Value = 3.2
    $ws.Cells.Item(105, 13).Value = 2.6
    $ws.Cells.Item(105, 14).Value = 2.05
    $ws.Cells.Item(105, 15).Value = 3.3
    $ws.Cells.Item(105, 16).Value = 3.25
    $ws.Cells.Item(105, 17).Value = -0.25
    $ws.Cells.Item(105, 18).Value = 1.85
    $ws.Cells.Item(105, 19).Value = 1.95
    $ws.Cells.Item(105, 20).Value = 2.75
    $ws.Cells.Item(105, 21).Value = 1.9
    $ws.Cells.Item(105, 22).Value = 1.9
    $ws.Cells.Item(105, 23).Value = 1.05
    $ws.Cells.Item(105, 24).Value = -1
    $ws.Cells.Item(105, 25).Value = -1
    $ws.Cells.Item(105, 26).Value = 0.8500000000000001
    $ws.Cells.Item(105, 27).Value = -1
    $ws.Cells.Item(105, 28).Value = 0.8999999999999999
    $ws.Cells.Item(105, 29).Value = -1
    # Row 106
    $ws.Cells.Item(106, 2).Value = 7630065
    $ws.Cells.Item(106, 6).Value = "Bray Wanderers"
    $ws.Cells.Item(106, 7).Value = "UCD"
    $ws.Cells.Item(106, 8).Value = 1
    $ws.Cells.Item(106, 9).Value = 2
    $ws.Cells.Item(106, 10).Value = "A"
    $ws.Cells.Item(106, 11).Value = 2.3
    $ws.Cells.Item(106, 12).Value = 3.4
    $ws.Cells.Item(106, 13).Value = 2.7
    $ws.Cells.Item(106, 14).Value = 1.8
    $ws.Cells.Item(106, 15).Value = 3.75
    $ws.Cells.Item(106, 16).Value = 3.8
    $ws.Cells.Item(106, 17).Value = -0.5
    $ws.Cells.Item(106, 18).Value = 1.825
    $ws.Cells.Item(106, 19).Value = 2.025
    $ws.Cells.Item(106, 20).Value = 2.75
    $ws.Cells.Item(106, 21).Value = 2.025
    $ws.Cells.Item(106, 22).Value = 1.825
    $ws.Cells.Item(106, 23).Value = -1
    $ws.Cells.Item(106, 24).Value = -1
    $ws.Cells.Item(106, 25).Value = 2.8
    $ws.Cells.Item(106, 26).Value = -1
    $ws.Cells.Item(106, 27).Value = 1.025
    $ws.Cells.Item(106, 28).Value = 0.5125
    $ws.Cells.Item(106, 29).Value = -0.5
    # Row 124
    $ws.Cells.Item(124, 2).Value = 7630078
    $ws.Cells.Item(124, 6).Value = "Cobh Ramblers"
    $ws.Cells.Item(124, 7).Value = "UCD"
    $ws.Cells.Item(124, 8).Value = 1
    $ws.Cells.Item(124, 9).Value = 2
    $ws.Cells.Item(124, 10).Value = "A"
    $ws.Cells.Item(124, 11).Value = 2.3
    $ws.Cells.Item(124, 12).Value = 3.2
    $ws.Cells.Item(124, 13).Value = 2.7
    $ws.Cells.Item(124, 14).Value = 2.45
    $ws.Cells.Item(124, 15).Value = 3.1
    $ws.Cells.Item(124, 16).Value = 2.625
    $ws.Cells.Item(124, 17).Value = 0
    $ws.Cells.Item(124, 18).Value = 1.8
    $ws.Cells.Item(124, 19).Value = 2.05
    $ws.Cells.Item(124, 20).Value = 2.25
    $ws.Cells.Item(124, 21).Value = 1.875
    $ws.Cells.Item(124, 22).Value = 1.975
    $ws.Cells.Item(124, 23).Value = -1
    $ws.Cells.Item(124, 24).Value = -1
    $ws.Cells.Item(124, 25).Value = 1.625
    $ws.Cells.Item(124, 26).Value = -1
    $ws.Cells.Item(124, 27).Value = 1.05
    $ws.Cells.Item(124, 28).Value = 0.875
    $ws.Cells.Item(124, 29).Value = -1
    # Row 126
    $ws.Cells.Item(126, 2).Value = 7630149
    $ws.Cells.Item(126, 6).Value = "Finn Harps"
    $ws.Cells.Item(126, 7).Value = "Kerry FC"
    $ws.Cells.Item(126, 8).Value = 0
    $ws.Cells.Item(126, 9).Value = 0
    $ws.Cells.Item(126, 10).Value = "D"
    $ws.Cells.Item(126, 11).Value = 1.615
    $ws.Cells.Item(126, 12).Value = 3.75
    $ws.Cells.Item(126, 13).Value = 4.5
    $ws.Cells.Item(126, 14).Value = 1.6
    $ws.Cells.Item(126, 15).Value = 3.75
    $ws.Cells.Item(126, 16).Value = 5
    $ws.Cells.Item(126, 17).Value = -0.75
    $ws.Cells.Item(126, 18).Value = 1.775
    $ws.Cells.Item(126, 19).Value = 2.025
    $ws.Cells.Item(126, 20).Value = 2.5
    $ws.Cells.Item(126, 21).Value = 1.95
    $ws.Cells.Item(126, 22).Value = 1.85
    $ws.Cells.Item(126, 23).Value = -1
    $ws.Cells.Item(126, 24).Value = 2.75
    $ws.Cells.Item(126, 25).Value = -1
    $ws.Cells.Item(126, 26).Value = -1
    $ws.Cells.Item(126, 27).Value = 1.025
    $ws.Cells.Item(126, 28).Value = -1
    $ws.Cells.Item(126, 29).Value = 0.8500000000000001
    # Row 127
    $ws.Cells.Item(127, 2).Value = 7630590
    $ws.Cells.Item(127, 6).Value = "Treaty United FC"
    $ws.Cells.Item(127, 7).Value = "Cork City"
    $ws.Cells.Item(127, 8).Value = 0
    $ws.Cells.Item(127, 9).Value = 0
    $ws.Cells.Item(127, 10).Value = "D"
    $ws.Cells.Item(127, 11).Value = 4.333
    $ws.Cells.Item(127, 12).Value = 3.7
    $ws.Cells.Item(127, 13).Value = 1.666
    $ws.Cells.Item(127, 14).Value = 4.2
    $ws.Cells.Item(127, 15).Value = 3.6
    $ws.Cells.Item(127, 16).Value = 1.727
    $ws.Cells.Item(127, 17).Value = 0.75
    $ws.Cells.Item(127, 18).Value = 1.8
    $ws.Cells.Item(127, 19).Value = 2
    $ws.Cells.Item(127, 20).Value = 2.5
    $ws.Cells.Item(127, 21).Value = 2
    $ws.Cells.Item(127, 22).Value = 1.8
    $ws.Cells.Item(127, 23).Value = -1
    $ws.Cells.Item(127, 24).Value = 2.6
    $ws.Cells.Item(127, 25).Value = -1
    $ws.Cells.Item(127, 26).Value = 0.8
    $ws.Cells.Item(127, 27).Value = -1
    $ws.Cells.Item(127, 28).Value = -1
    $ws.Cells.Item(127, 29).Value = 0.8
    # Row 129
    $ws.Cells.Item(129, 2).Value = 7630080
    $ws.Cells.Item(129, 6).Value = "Longford Town"
    $ws.Cells.Item(129, 7).Value = "UCD"
    $ws.Cells.Item(129, 8).Value = 1
    $ws.Cells.Item(129, 9).Value = 2
    $ws.Cells.Item(129, 10).Value = "A"
    $ws.Cells.Item(129, 11).Value = 2.375
    $ws.Cells.Item(129, 12).Value = 3.4
    $ws.Cells.Item(129, 13).Value = 2.625
    $ws.Cells.Item(129, 14).Value = 2.55
    $ws.Cells.Item(129, 15).Value = 3.1
    $ws.Cells.Item(129, 16).Value = 2.6
    $ws.Cells.Item(129, 17).Value = 0
    $ws.Cells.Item(129, 18).Value = 1.85
    $ws.Cells.Item(129, 19).Value = 1.95
    $ws.Cells.Item(129, 20).Value = 2.25
    $ws.Cells.Item(129, 21).Value = 2.025
    $ws.Cells.Item(129, 22).Value = 1.775
    $ws.Cells.Item(129, 23).Value = -1
    $ws.Cells.Item(129, 24).Value = -1
    $ws.Cells.Item(129, 25).Value = 1.6
    $ws.Cells.Item(129, 26).Value = -1
    $ws.Cells.Item(129, 27).Value = 0.95
    $ws.Cells.Item(129, 28).Value = 1.025
    $ws.Cells.Item(129, 29).Value = -1
    # Row 130
    $ws.Cells.Item(130, 2).Value = 7630150
    $ws.Cells.Item(130, 6).Value = "Treaty United FC"
    $ws.Cells.Item(130, 7).Value = "Finn Harps"
    $ws.Cells.Item(130, 8).Value = 0
    $ws.Cells.Item(130, 9).Value = 1
    $ws.Cells.Item(130, 10).Value = "A"
    $ws.Cells.Item(130, 11).Value = 1.95
    $ws.Cells.Item(130, 12).Value = 3.5
    $ws.Cells.Item(130, 13).Value = 3.3
    $ws.Cells.Item(130, 14).Value = 2.25
    $ws.Cells.Item(130, 15).Value = 3.25
    $ws.Cells.Item(130, 16).Value = 2.8
    $ws.Cells.Item(130, 17).Value = -0.25
    $ws.Cells.Item(130, 18).Value = 2.1
    $ws.Cells.Item(130, 19).Value = 1.775
    $ws.Cells.Item(130, 20).Value = 2.5
    $ws.Cells.Item(130, 21).Value = 1.975
    $ws.Cells.Item(130, 22).Value = 1.875
    $ws.Cells.Item(130, 23).Value = -1
    $ws.Cells.Item(130, 24).Value = -1
    $ws.Cells.Item(130, 25).Value = 1.8
    $ws.Cells.Item(130, 26).Value = -1
    $ws.Cells.Item(130, 27).Value = 0.7749999999999999
    $ws.Cells.Item(130, 28).Value = -1
    $ws.Cells.Item(130, 29).Value = 0.875
    # Row 131
    $ws.Cells.Item(131, 2).Value = 7630591
    $ws.Cells.Item(131, 6).Value = "Cork City"
    $ws.Cells.Item(131, 7).Value = "Athlone Town"
    $ws.Cells.Item(131, 8).Value = 1
    $ws.Cells.Item(131, 9).Value = 0
    $ws.Cells.Item(131, 10).Value = "H"
    $ws.Cells.Item(131, 11).Value = 1.4
    $ws.Cells.Item(131, 12).Value = 4.2
    $ws.Cells.Item(131, 13).Value = 6.5
    $ws.Cells.Item(131, 14).Value = 1.45
    $ws.Cells.Item(131, 15).Value = 4
    $ws.Cells.Item(131, 16).Value = 6.5
    $ws.Cells.Item(131, 17).Value = -1
    $ws.Cells.Item(131, 18).Value = 1.8
    $ws.Cells.Item(131, 19).Value = 2
    $ws.Cells.Item(131, 20).Value = 2.5
    $ws.Cells.Item(131, 21).Value = 1.9
    $ws.Cells.Item(131, 22).Value = 1.9
    $ws.Cells.Item(131, 23).Value = 0.45
    $ws.Cells.Item(131, 24).Value = -1
    $ws.Cells.Item(131, 25).Value = -1
    $ws.Cells.Item(131, 26).Value = 0
    $ws.Cells.Item(131, 27).Value = -0
    $ws.Cells.Item(131, 28).Value = -1
    $ws.Cells.Item(131, 29).Value = 0.8999999999999999
    # Row 132
    $ws.Cells.Item(132, 2).Value = 7630592
    $ws.Cells.Item(132, 6).Value = "Kerry FC"
    $ws.Cells.Item(132, 7).Value = "Wexford FC"
    $ws.Cells.Item(132, 8).Value = 1
    $ws.Cells.Item(132, 9).Value = 1
    $ws.Cells.Item(132, 10).Value = "D"
    $ws.Cells.Item(132, 11).Value = 4
    $ws.Cells.Item(132, 12).Value = 3.6
    $ws.Cells.Item(132, 13).Value = 1.75
    $ws.Cells.Item(132, 14).Value = 4.5
    $ws.Cells.Item(132, 15).Value = 3.75
    $ws.Cells.Item(132, 16).Value = 1.65
    $ws.Cells.Item(132, 17).Value = 0.75
    $ws.Cells.Item(132, 18).Value = 2
    $ws.Cells.Item(132, 19).Value = 1.8
    $ws.Cells.Item(132, 20).Value = 2.5
    $ws.Cells.Item(132, 21).Value = 1.85
    $ws.Cells.Item(132, 22).Value = 1.95
    $ws.Cells.Item(132, 23).Value = -1
    $ws.Cells.Item(132, 24).Value = 2.75
    $ws.Cells.Item(132, 25).Value = -1
    $ws.Cells.Item(132, 26).Value = 1
    $ws.Cells.Item(132, 27).Value = -1
    $ws.Cells.Item(132, 28).Value = -1
    $ws.Cells.Item(132, 29).Value = 0.95
    # Row 133
    $ws.Cells.Item(133, 2).Value = 7630079
    $ws.Cells.Item(133, 6).Value = "Bray Wanderers"
    $ws.Cells.Item(133, 7).Value = "Cobh Ramblers"
    $ws.Cells.Item(133, 8).Value = 2
    $ws.Cells.Item(133, 9).Value = 2
    $ws.Cells.Item(133, 10).Value = "D"
    $ws.Cells.Item(133, 11).Value = 1.833
    $ws.Cells.Item(133, 12).Value = 3.5
    $ws.Cells.Item(133, 13).Value = 3.6
    $ws.Cells.Item(133, 14).Value = 2.4
    $ws.Cells.Item(133, 15).Value = 3.25
    $ws.Cells.Item(133, 16).Value = 2.6
    $ws.Cells.Item(133, 17).Value = 0
    $ws.Cells.Item(133, 18).Value = 1.8
    $ws.Cells.Item(133, 19).Value = 2
    $ws.Cells.Item(133, 20).Value = 2.75
    $ws.Cells.Item(133, 21).Value = 1.975
    $ws.Cells.Item(133, 22).Value = 1.825
    $ws.Cells.Item(133, 23).Value = -1
    $ws.Cells.Item(133, 24).Value = 2.25
    $ws.Cells.Item(133, 25).Value = -1
    $ws.Cells.Item(133, 26).Value = 0
    $ws.Cells.Item(133, 27).Value = -0
    $ws.Cells.Item(133, 28).Value = 0.9750000000000001
    $ws.Cells.Item(133, 29).Value = -1
    # Row 136
    $ws.Cells.Item(136, 2).Value = 7630081
    $ws.Cells.Item(136, 6).Value = "Finn Harps"
    $ws.Cells.Item(136, 7).Value = "Bray Wanderers"
    $ws.Cells.Item(136, 8).Value = 2
    $ws.Cells.Item(136, 9).Value = 0
    $ws.Cells.Item(136, 10).Value = "H"
    $ws.Cells.Item(136, 11).Value = 2.4
    $ws.Cells.Item(136, 12).Value = 3.4
    $ws.Cells.Item(136, 13).Value = 2.6
    $ws.Cells.Item(136, 14).Value = 2.1
    $ws.Cells.Item(136, 15).Value = 3.4
    $ws.Cells.Item(136, 16).Value = 3
    $ws.Cells.Item(136, 17).Value = -0.25
    $ws.Cells.Item(136, 18).Value = 1.85
    $ws.Cells.Item(136, 19).Value = 1.95
    $ws.Cells.Item(136, 20).Value = 2.5
    $ws.Cells.Item(136, 21).Value = 1.85
    $ws.Cells.Item(136, 22).Value = 1.95
    $ws.Cells.Item(136, 23).Value = 1.1
    $ws.Cells.Item(136, 24).Value = -1
    $ws.Cells.Item(136, 25).Value = -1
    $ws.Cells.Item(136, 26).Value = 0.8500000000000001
    $ws.Cells.Item(136, 27).Value = -1
    $ws.Cells.Item(136, 28).Value = -1
    $ws.Cells.Item(136, 29).Value = 0.95
    # Row 137
    $ws.Cells.Item(137, 2).Value = 7630593
    $ws.Cells.Item(137, 6).Value = "Athlone Town"
    $ws.Cells.Item(137, 7).Value = "Treaty United FC"
    $ws.Cells.Item(137, 8).Value = 1
    $ws.Cells.Item(137, 9).Value = 0
    $ws.Cells.Item(137, 10).Value = "H"
    $ws.Cells.Item(137, 11).Value = 2.9
    $ws.Cells.Item(137, 12).Value = 3.3
    $ws.Cells.Item(137, 13).Value = 2.25
    $ws.Cells.Item(137, 14).Value = 2.5
    $ws.Cells.Item(137, 15).Value = 3.3
    $ws.Cells.Item(137, 16).Value = 2.6
    $ws.Cells.Item(137, 17).Value = 0
    $ws.Cells.Item(137, 18).Value = 1.85
    $ws.Cells.Item(137, 19).Value = 1.95
    $ws.Cells.Item(137, 20).Value = 2.75
    $ws.Cells.Item(137, 21).Value = 2
    $ws.Cells.Item(137, 22).Value = 1.8
    $ws.Cells.Item(137, 23).Value = 1.5
    $ws.Cells.Item(137, 24).Value = -1
    $ws.Cells.Item(137, 25).Value = -1
    $ws.Cells.Item(137, 26).Value = 0.8500000000000001
    $ws.Cells.Item(137, 27).Value = -1
    $ws.Cells.Item(137, 28).Value = -1
    $ws.Cells.Item(137, 29).Value = 0.8
    # Row 138
    $ws.Cells.Item(138, 2).Value = 7630082
    $ws.Cells.Item(138, 6).Value = "UCD"
    $ws.Cells.Item(138, 7).Value = "Kerry FC"
    $ws.Cells.Item(138, 8).Value = 1
    $ws.Cells.Item(138, 9).Value = 1
    $ws.Cells.Item(138, 10).Value = "D"
    $ws.Cells.Item(138, 11).Value = 1.533
    $ws.Cells.Item(138, 12).Value = 4
    $ws.Cells.Item(138, 13).Value = 5.25
    $ws.Cells.Item(138, 14).Value = 1.75
    $ws.Cells.Item(138, 15).Value = 3.5
    $ws.Cells.Item(138, 16).Value = 4.2
    $ws.Cells.Item(138, 17).Value = -0.5
    $ws.Cells.Item(138, 18).Value = 1.775
    $ws.Cells.Item(138, 19).Value = 2.025
    $ws.Cells.Item(138, 20).Value = 2.5
    $ws.Cells.Item(138, 21).Value = 1.95
    $ws.Cells.Item(138, 22).Value = 1.85
    $ws.Cells.Item(138, 23).Value = -1
    $ws.Cells.Item(138, 24).Value = 2.5
    $ws.Cells.Item(138, 25).Value = -1
    $ws.Cells.Item(138, 26).Value = -1
    $ws.Cells.Item(138, 27).Value = 1.025
    $ws.Cells.Item(138, 28).Value = -1
    $ws.Cells.Item(138, 29).Value = 0.8500000000000001

Write-Host "Row rotation update applied."
